$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'61.350.17"
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = "'  -4.64%  "
$ws.Range('E2').Style = 'Normal'
$ws.Range('D3').Value = "'3.331.65"
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = "'  -4.81%  "
$ws.Range('E3').Style = 'Normal'
$ws.Range('D4').Value = "'1.00"
$ws.Range('D4').Style = 'Normal'
$ws.Range('D5').Value = "'569.56"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = "'  -3.43%  "
$ws.Range('E5').Style = 'Normal'
$ws.Range('D6').Value = "'128.30"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = "'  -4.34%  "
$ws.Range('E6').Style = 'Normal'
$ws.Range('E7').Value = "'  +0.02%  "
$ws.Range('E7').Style = 'Normal'
$ws.Range('D8').Value = "'3.329.91"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = "'  -4.85%  "
$ws.Range('E8').Style = 'Normal'
$ws.Range('D9').Value = "'0.477"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = "'  -1.76%  "
$ws.Range('E9').Style = 'Normal'
$ws.Range('D10').Value = "'7.39"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = "'  -4.19%  "
$ws.Range('E10').Style = 'Normal'
$ws.Range('D11').Value = "'0.118"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = "'  -5.05%  "
$ws.Range('E11').Style = 'Normal'
$ws.Range('D12').Value = "'0.375"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = "'  -3.02%  "
$ws.Range('E12').Style = 'Normal'
$ws.Range('D13').Value = "'3.897.43"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = "'  -4.81%  "
$ws.Range('E13').Style = 'Normal'
$ws.Range('D14').Value = "'0.120"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = "'  -0.44%  "
$ws.Range('E14').Style = 'Normal'
$ws.Range('D15').Value = "'3.317.65"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = "'  -5.17%  "
$ws.Range('E15').Style = 'Normal'
$ws.Range('D16').Value = "'0.0000168"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = "'  -6.32%  "
$ws.Range('E16').Style = 'Normal'
$ws.Range('B17').Value = "'Avalanche"
$ws.Range('B17').Style = 'Normal'
$ws.Range('C17').Value = "'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range('C17').Style = 'Normal'
$ws.Range('D17').Value = "'24.90"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = "'  -1.46%  "
$ws.Range('E17').Style = 'Normal'
$ws.Range('B18').Value = "'WrappedBTC"
$ws.Range('B18').Style = 'Normal'
$ws.Range('C18').Value = "'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range('C18').Style = 'Normal'
$ws.Range('D18').Value = "'61.452.79"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = "'  -4.39%  "
$ws.Range('E18').Style = 'Normal'
$ws.Range('D19').Value = "'5.63"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = "'  -2.29%  "
$ws.Range('E19').Style = 'Normal'
$ws.Range('D20').Value = "'13.30"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = "'  -1.84%  "
$ws.Range('E20').Style = 'Normal'
$ws.Range('D21').Value = "'9.04"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = "'  -10.04%  "
$ws.Range('E21').Style = 'Normal'
$ws.Range('D22').Value = "'356.56"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = "'  -7.81%  "
$ws.Range('E22').Style = 'Normal'
$ws.Range('D23').Value = "'0.555"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = "'  -4.30%  "
$ws.Range('E23').Style = 'Normal'
$ws.Range('E24').Value = "'  -0.06%  "
$ws.Range('E24').Style = 'Normal'
$ws.Range('D25').Value = "'3.461.51"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = "'  -4.87%  "
$ws.Range('E25').Style = 'Normal'
$ws.Range('D26').Value = "'69.85"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = "'  -5.91%  "
$ws.Range('E26').Style = 'Normal'
$ws.Range('D27').Value = "'0.0000108"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = "'  -6.75%  "
$ws.Range('E27').Style = 'Normal'
$ws.Range('D28').Value = "'0.998"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = "'  +0.01%  "
$ws.Range('E28').Style = 'Normal'
$ws.Range('D29').Value = "'7.22"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = "'  -0.83%  "
$ws.Range('E29').Style = 'Normal'
$ws.Range('E30').Value = "'  -2.87%  "
$ws.Range('E30').Style = 'Normal'
$ws.Range('D31').Value = "'7.93"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = "'  -2.26%  "
$ws.Range('E31').Style = 'Normal'
$ws.Range('D32').Value = "'2.12"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = "'  -6.04%  "
$ws.Range('E32').Style = 'Normal'
$ws.Range('E33').Value = "'  -0.06%  "
$ws.Range('E33').Style = 'Normal'
$ws.Range('D34').Value = "'0.149"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = "'  -3.86%  "
$ws.Range('E34').Style = 'Normal'
$ws.Range('D35').Value = "'3.360.05"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = "'  -4.74%  "
$ws.Range('E35').Style = 'Normal'
$ws.Range('D36').Value = "'22.55"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = "'  -3.28%  "
$ws.Range('E36').Style = 'Normal'
$ws.Range('E37').Value = "'  -1.06%  "
$ws.Range('E37').Style = 'Normal'
$ws.Range('D38').Value = "'6.77"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = "'  -2.02%  "
$ws.Range('E38').Style = 'Normal'
$ws.Range('D39').Value = "'161.64"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = "'  -2.01%  "
$ws.Range('E39').Style = 'Normal'
$ws.Range('D40').Value = "'1.49"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = "'  -3.13%  "
$ws.Range('E40').Style = 'Normal'
$ws.Range('D41').Value = "'0.0758"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = "'  -3.44%  "
$ws.Range('E41').Style = 'Normal'
$ws.Range('D42').Value = "'1.00"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = "'  +0.00%  "
$ws.Range('E42').Style = 'Normal'
$ws.Range('D43').Value = "'41.31"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = "'  -1.74%  "
$ws.Range('E43').Style = 'Normal'
$ws.Range('B44').Value = "'Filecoin"
$ws.Range('B44').Style = 'Normal'
$ws.Range('C44').Value = "'https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range('C44').Style = 'Normal'
$ws.Range('D44').Value = "'4.36"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = "'  -1.10%  "
$ws.Range('E44').Style = 'Normal'
$ws.Range('B45').Value = "'Mantle"
$ws.Range('B45').Style = 'Normal'
$ws.Range('C45').Value = "'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range('C45').Style = 'Normal'
$ws.Range('D45').Value = "'0.748"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = "'  -7.42%  "
$ws.Range('E45').Style = 'Normal'
$ws.Range('D46').Value = "'1.13"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = "'  -3.93%  "
$ws.Range('E46').Style = 'Normal'
$ws.Range('D47').Value = "'1.58"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = "'  -4.35%  "
$ws.Range('E47').Style = 'Normal'
$ws.Range('D48').Value = "'22.33"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = "'  -8.56%  "
$ws.Range('E48').Style = 'Normal'
$ws.Range('D49').Value = "'6.69"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = "'  -1.75%  "
$ws.Range('E49').Style = 'Normal'
$ws.Range('D50').Value = "'0.871"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = "'  -5.19%  "
$ws.Range('E50').Style = 'Normal'
$ws.Range('B51').Value = "'InjectiveProtocol"
$ws.Range('B51').Style = 'Normal'
$ws.Range('C51').Value = "'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range('C51').Style = 'Normal'
$ws.Range('D51').Value = "'21.24"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = "'  +1.20%  "
$ws.Range('E51').Style = 'Normal'
